# Edit: "Added modulo function to git code."
#   - Slide with the "Group Exercise!" tasks list: the bullet that used to
#     read "Eliot: increment" should read "Eliot: mod" instead.
#
# We walk the deck looking for the run of text rather than hard-coding a
# slide/shape number, then do a surgical in-place replacement of just that
# substring via TextRange.Characters(start, length) so every other run /
# paragraph in the text body (bullets, fonts, bullet numbering, etc.) is
# left completely untouched.

$p = $ppt.ActivePresentation

$oldText = "Eliot: increment"
$newText = "Eliot: mod"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -eq $false) {
            continue
        }
        $textRange = $shape.TextFrame.TextRange
        $fullText = $textRange.Text
        $idx = $fullText.IndexOf($oldText)
        if ($idx -ge 0) {
            $hit = $textRange.Characters($idx + 1, $oldText.Length)
            $hit.Text = $newText
        }
    }
}
